$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" conversion note text (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$old1 = "$([char]0x2705) 1000 Bs = 14.6 = 59868.76 pesos"
$new1 = "$([char]0x2705) 1000 Bs = 14.62 = 60032.47 pesos"
$old2 = "$([char]0x2705) 59868.76 pesos = 14.55 = 971.9 Bs"
$new2 = "$([char]0x2705) 60032.47 pesos = 14.59 = 977.54 Bs"

$text = $wsHoja1.Range("A1").Text
$text = $text.Replace($old1, $new1)
$text = $text.Replace($old2, $new2)
$wsHoja1.Range("A1").Value = $text

# --- Update the "tasas" sheet rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 68.38
$wsTasas.Range("O10").Value = 4105.02
$wsTasas.Range("N12").Value = 4114.9
$wsTasas.Range("O12").Value = 67.005
